$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.971.74'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.20%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.860.82'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.52%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = "Normal"

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.97'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.50%  '

# Row 6
$ws.Range('E6').Value = '  -0.08%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5129'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.84%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3820'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.53%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08267'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.06%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.109'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.42%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.51'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.07%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.188'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.90%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.48'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.73%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.850.59'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.77%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.271'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.41%  '

# Row 16
$ws.Range('E16').Value = '  -0.10%  '

# Row 17
$ws.Range('E17').Value = '  -0.41%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '90.36'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.55%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06642'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.08%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.65'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.98%  '

# Row 21
$ws.Range('E21').Value = '  -0.08%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.015'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.00%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.009.58'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.17%  '

# Row 24
$ws.Range('E24').Value = '  -2.94%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.244'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.01%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.068.43'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.35%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.500'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.62%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '157.24'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.16%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '20.46'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.09%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.36'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.28%  '

# Row 31
$ws.Range('E31').Value = '  +1.14%  '

# Row 32
$ws.Range('E32').Value = '  -3.08%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.827'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +4.14%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.590'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.18%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.410'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.92%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02406'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.14%  '

# Row 37
$ws.Range('E37').Value = '  -1.14%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2192'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.11%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6523'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.63%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.195'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.84%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.981'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.15%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.209'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.29%  '

# Row 43
$ws.Range('E43').Value = '  -3.03%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6102'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.09%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '12.97'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.18%  '

# Row 46
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('B46').Style = "Normal"
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C46').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.276'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.38%  '

# Row 47
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('B47').Style = "Normal"
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C47').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.670'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.11%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.012'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.43%  '

# Row 49
$ws.Range('E49').Value = '  -1.33%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '120.61'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.54%  '

# Row 51
$ws.Range('E51').Value = '  -3.25%  '
